$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H53").Value = 1934.3846
$ws_ALC.Range("I53").Value = 249.8
$ws_ALC.Range("K53").Value = 249.8
$ws_ALC.Range("M53").Value = 387.2
$ws_ALC.Range("H64").Value = 9266.333000000001
$ws_ALC.Range("J64").Value = 9733.833000000001
$ws_ALC.Range("L64").Value = 9733.833000000001
$ws_ALC.Range("N64").Value = -10229.833
$ws_ALC.Range("H67").Value = 9266.333000000001
$ws_ALC.Range("J67").Value = 9733.833000000001
$ws_ALC.Range("L67").Value = 9733.833000000001
$ws_ALC.Range("N67").Value = -11449.833
$ws_ALC.Range("H121").Value = 2800
$ws_ALC.Range("J121").Value = 2800
$ws_ALC.Range("L121").Value = 8400
$ws_ALC.Range("N121").Value = -11894
$ws_ALC.Range("H132").Value = 1525.3265
$ws_ALC.Range("I132").Value = 1535.659
$ws_ALC.Range("J132").Value = 1434.4
$ws_ALC.Range("K132").Value = 4606.977000000001
$ws_ALC.Range("L132").Value = 4303.200000000001
$ws_ALC.Range("M132").Value = -2076.977000000001
$ws_ALC.Range("N132").Value = -9363.200000000001
$ws_ALC.Range("H135").Value = 1059.9166
$ws_ALC.Range("I135").Value = 816.4737
$ws_ALC.Range("J135").Value = 1985
$ws_ALC.Range("K135").Value = 7348.263300000001
$ws_ALC.Range("L135").Value = 17865
$ws_ALC.Range("M135").Value = -4813.263300000001
$ws_ALC.Range("N135").Value = -22935
$ws_ALC.Range("H138").Value = 3329.027
$ws_ALC.Range("J138").Value = 3077.1614
$ws_ALC.Range("L138").Value = 9231.484199999999
$ws_ALC.Range("N138").Value = -19511.4842

# Sheet: ARM
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H2").Value = 6929.6665
$ws_ARM.Range("I2").Value = 728.9091
$ws_ARM.Range("K2").Value = 728.9091
$ws_ARM.Range("M2").Value = -615.9091
$ws_ARM.Range("H32").Value = 16131328
$ws_ARM.Range("I32").Value = 16395301
$ws_ARM.Range("K32").Value = 16395301
$ws_ARM.Range("M32").Value = -16395014
$ws_ARM.Range("H45").Value = 3241.7222
$ws_ARM.Range("I45").Value = 2177.0667
$ws_ARM.Range("K45").Value = 2177.0667
$ws_ARM.Range("M45").Value = -1800.0667
$ws_ARM.Range("H74").Value = 13890682
$ws_ARM.Range("I74").Value = 14494368
$ws_ARM.Range("J74").Value = 5900
$ws_ARM.Range("K74").Value = 14494368
$ws_ARM.Range("L74").Value = 5900
$ws_ARM.Range("M74").Value = -14493494
$ws_ARM.Range("N74").Value = -7648
$ws_ARM.Range("H77").Value = 13890682
$ws_ARM.Range("I77").Value = 14494368
$ws_ARM.Range("J77").Value = 5900
$ws_ARM.Range("K77").Value = 72471840
$ws_ARM.Range("L77").Value = 29500
$ws_ARM.Range("M77").Value = -72467472
$ws_ARM.Range("N77").Value = -38236
$ws_ARM.Range("H116").Value = 6929.6665
$ws_ARM.Range("I116").Value = 728.9091
$ws_ARM.Range("K116").Value = 728.9091
$ws_ARM.Range("M116").Value = 1565.0909

# Sheet: BSM
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H3").Value = 6929.6665
$ws_BSM.Range("I3").Value = 728.9091
$ws_BSM.Range("K3").Value = 728.9091
$ws_BSM.Range("M3").Value = -614.9091
$ws_BSM.Range("H22").Value = 2255.4285
$ws_BSM.Range("J22").Value = 3374.5
$ws_BSM.Range("L22").Value = 3374.5
$ws_BSM.Range("N22").Value = -3720.5
$ws_BSM.Range("H86").Value = 2621.5715
$ws_BSM.Range("J86").Value = 3197.3333
$ws_BSM.Range("L86").Value = 3197.3333
$ws_BSM.Range("N86").Value = -5443.3333
$ws_BSM.Range("H89").Value = 2621.5715
$ws_BSM.Range("J89").Value = 3197.3333
$ws_BSM.Range("L89").Value = 15986.6665
$ws_BSM.Range("N89").Value = -27218.6665
$ws_BSM.Range("H99").Value = 2030.091
$ws_BSM.Range("I99").Value = 2001.125
$ws_BSM.Range("K99").Value = 2001.125
$ws_BSM.Range("M99").Value = -503.125

# Sheet: CRP
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 77863.14
$ws_CRP.Range("I31").Value = 3632.125
$ws_CRP.Range("K31").Value = 3632.125
$ws_CRP.Range("M31").Value = -3337.125
$ws_CRP.Range("H34").Value = 77863.14
$ws_CRP.Range("I34").Value = 3632.125
$ws_CRP.Range("K34").Value = 3632.125
$ws_CRP.Range("M34").Value = -3430.125
$ws_CRP.Range("H58").Value = 4431.8335
$ws_CRP.Range("I58").Value = 2044.5333
$ws_CRP.Range("J58").Value = 8410.666999999999
$ws_CRP.Range("K58").Value = 2044.5333
$ws_CRP.Range("L58").Value = 8410.666999999999
$ws_CRP.Range("M58").Value = -1841.5333
$ws_CRP.Range("N58").Value = -8816.666999999999
$ws_CRP.Range("H99").Value = 3644.5
$ws_CRP.Range("I99").Value = 3644.5
$ws_CRP.Range("K99").Value = 3644.5
$ws_CRP.Range("M99").Value = -2146.5
$ws_CRP.Range("H118").Value = 60000
$ws_CRP.Range("J118").Value = 60000
$ws_CRP.Range("L118").Value = 60000
$ws_CRP.Range("N118").Value = -63314
$ws_CRP.Range("H126").Value = 3644.5
$ws_CRP.Range("I126").Value = 3644.5
$ws_CRP.Range("K126").Value = 10933.5
$ws_CRP.Range("M126").Value = -8463.5
$ws_CRP.Range("H132").Value = 5885.0713
$ws_CRP.Range("I132").Value = 3262.7
$ws_CRP.Range("K132").Value = 9788.099999999999
$ws_CRP.Range("M132").Value = -7258.099999999999
$ws_CRP.Range("H134").Value = 4965.231
$ws_CRP.Range("I134").Value = 2775
$ws_CRP.Range("K134").Value = 8325
$ws_CRP.Range("M134").Value = -5790
$ws_CRP.Range("H136").Value = 4431.8335
$ws_CRP.Range("I136").Value = 2044.5333
$ws_CRP.Range("J136").Value = 8410.666999999999
$ws_CRP.Range("K136").Value = 6133.5999
$ws_CRP.Range("L136").Value = 25232.001
$ws_CRP.Range("M136").Value = -3583.5999
$ws_CRP.Range("N136").Value = -30332.001

# Sheet: CUL
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H5").Value = 15987.8
$ws_CUL.Range("J5").Value = 30204.8
$ws_CUL.Range("L5").Value = 90614.39999999999
$ws_CUL.Range("N5").Value = -90838.39999999999
$ws_CUL.Range("H55").Value = 1139.3846
$ws_CUL.Range("I55").Value = 857
$ws_CUL.Range("J55").Value = 1774.75
$ws_CUL.Range("K55").Value = 2571
$ws_CUL.Range("L55").Value = 5324.25
$ws_CUL.Range("M55").Value = -2394
$ws_CUL.Range("N55").Value = -5678.25
$ws_CUL.Range("H97").Value = 373.0909
$ws_CUL.Range("I97").Value = 232.88889
$ws_CUL.Range("K97").Value = 698.6666700000001
$ws_CUL.Range("M97").Value = -202.6666700000001
$ws_CUL.Range("H122").Value = 6482432.5
$ws_CUL.Range("I122").Value = 6410763.5
$ws_CUL.Range("K122").Value = 57696871.5
$ws_CUL.Range("M122").Value = -57694421.5
$ws_CUL.Range("H135").Value = 15987.8
$ws_CUL.Range("J135").Value = 30204.8
$ws_CUL.Range("L135").Value = 271843.2
$ws_CUL.Range("N135").Value = -276913.2
$ws_CUL.Range("H139").Value = 4257.0835
$ws_CUL.Range("I139").Value = 2487.2
$ws_CUL.Range("K139").Value = 7461.599999999999
$ws_CUL.Range("M139").Value = -2321.599999999999

# Sheet: LTW
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H82").Value = 4346.75
$ws_LTW.Range("J82").Value = 3619.0667
$ws_LTW.Range("L82").Value = 3619.0667
$ws_LTW.Range("N82").Value = -4341.066699999999
$ws_LTW.Range("H85").Value = 4346.75
$ws_LTW.Range("J85").Value = 3619.0667
$ws_LTW.Range("L85").Value = 3619.0667
$ws_LTW.Range("N85").Value = -6115.066699999999
$ws_LTW.Range("H132").Value = 11792.083
$ws_LTW.Range("I132").Value = 2666.6667
$ws_LTW.Range("J132").Value = 14833.889
$ws_LTW.Range("K132").Value = 8000.000100000001
$ws_LTW.Range("L132").Value = 44501.667
$ws_LTW.Range("M132").Value = -5470.000100000001
$ws_LTW.Range("N132").Value = -49561.667
$ws_LTW.Range("H136").Value = 4307.5884
$ws_LTW.Range("J136").Value = 11414.571
$ws_LTW.Range("L136").Value = 34243.713
$ws_LTW.Range("N136").Value = -39343.713

# Sheet: WVR
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H81").Value = 2535.3044
$ws_WVR.Range("J81").Value = 9834
$ws_WVR.Range("L81").Value = 19668
$ws_WVR.Range("N81").Value = -21790
$ws_WVR.Range("H84").Value = 2535.3044
$ws_WVR.Range("J84").Value = 9834
$ws_WVR.Range("L84").Value = 98340
$ws_WVR.Range("N84").Value = -108948
$ws_WVR.Range("H107").Value = 602.93335
$ws_WVR.Range("I107").Value = 449.42856
$ws_WVR.Range("J107").Value = 737.25
$ws_WVR.Range("K107").Value = 1348.28568
$ws_WVR.Range("L107").Value = 2211.75
$ws_WVR.Range("M107").Value = 571.71432
$ws_WVR.Range("N107").Value = -6051.75
$ws_WVR.Range("H138").Value = 70429
$ws_WVR.Range("J138").Value = 70429
$ws_WVR.Range("L138").Value = 70429
$ws_WVR.Range("N138").Value = -80709
